$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.722.71'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '2.641.69'
$ws.Range('E3').Value = '  +1.26%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.22'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.83'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.600'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.25%  '
$ws.Range('E9').Value = '  +0.47%  '
$ws.Range('E10').Value = '  +0.08%  '
$ws.Range('E11').Value = '  +2.15%  '
$ws.Range('E12').Value = '  -0.83%  '
$ws.Range('D13').Value = '3.110.15'
$ws.Range('E13').Value = '  +1.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.11'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +11.45%  '
$ws.Range('D15').Value = '60.687.17'
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').Value = '2.655.77'
$ws.Range('E18').Value = '  +2.10%  '
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '349.72'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('E21').Value = '  -1.11%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('E23').Value = '  +1.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.92'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.12%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  +0.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.19'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.24%  '
$ws.Range('E28').Value = '  +9.32%  '
$ws.Range('E29').Value = '  +0.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.80'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +6.82%  '
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('E33').Value = '  +1.46%  '
$ws.Range('E34').Value = '  +7.63%  '
$ws.Range('E35').Value = '  +3.65%  '
$ws.Range('E36').Value = '  +6.58%  '
$ws.Range('E37').Value = '  +2.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '339.05'
$ws.Range('D38').ClearFormats()
$ws.Range('E39').Value = '  +4.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.905'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +6.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '38.32'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.17'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.25%  '
$ws.Range('E43').Value = '  +2.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.24'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.39%  '
$ws.Range('E45').Value = '  +2.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0562'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '132.87'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.54'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('E50').Value = '  +0.25%  '
$ws.Range('D51').Value = '2.085.89'
$ws.Range('E51').Value = '  +2.14%  '
